$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '41.366.96'
$ws.Range('E2').Value = '  -1.57%  '

# Row 3
$ws.Range('D3').Value = '2.194.25'
$ws.Range('E3').Value = '  -1.37%  '

# Row 4
$ws.Range('E4').Value = '  +0.17%  '

# Row 5
$ws.Range('D5').Value = '''252.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.42%  '

# Row 6
$ws.Range('D6').Value = '''0.624'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.19%  '

# Row 7
$ws.Range('D7').Value = '''68.85'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.69%  '

# Row 8
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').Value = '''0.586'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.86%  '

# Row 10
$ws.Range('D10').Value = '''38.12'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.58%  '

# Row 11
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '''58.09'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.15%  '

# Row 12
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').Value = '''0.0945'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.25%  '

# Row 13
$ws.Range('D13').Value = '''7.12'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.07%  '

# Row 14
$ws.Range('D14').Value = '''0.105'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.68%  '

# Row 15
$ws.Range('D15').Value = '2.534.00'
$ws.Range('E15').Value = '  -0.85%  '

# Row 16
$ws.Range('D16').Value = '''14.70'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.53%  '

# Row 17
$ws.Range('D17').Value = '''0.875'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.13%  '

# Row 18
$ws.Range('D18').Value = '2.185.41'
$ws.Range('E18').Value = '  -2.16%  '

# Row 19
$ws.Range('D19').Value = '41.336.96'
$ws.Range('E19').Value = '  -1.51%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0948'
$ws.Range('E20').Value = '  -0.56%  '

# Row 21
$ws.Range('D21').Value = '''6.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.36%  '

# Row 22
$ws.Range('D22').Value = '''71.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.93%  '

# Row 23
$ws.Range('D23').Value = '''231.53'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.33%  '

# Row 24
$ws.Range('D24').Value = '''2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.42%  '

# Row 25
$ws.Range('D25').Value = '''11.81'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +19.51%  '

# Row 26
$ws.Range('D26').Value = '''3.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.14%  '

# Row 27
$ws.Range('E27').Value = '  -0.03%  '

# Row 28
$ws.Range('D28').Value = '''2.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.26%  '

# Row 29
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').Value = '''3.76'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.19%  '

# Row 30
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''2.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.40%  '

# Row 31
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '''170.41'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.34%  '

# Row 32
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '''20.55'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.79%  '

# Row 33
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = '''0.121'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.58%  '

# Row 34
$ws.Range('E34').Value = '  -1.64%  '

# Row 35
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '''5.51'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.18%  '

# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.0729'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.32%  '

# Row 37
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').Value = '''25.96'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +15.91%  '

# Row 38
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').Value = '''4.61'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.87%  '

# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''3.99'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.33%  '

# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.0299'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.56%  '

# Row 41
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = '''2.22'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.15%  '

# Row 42
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').Value = '''5.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.49%  '

# Row 43
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').Value = '''11.95'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +17.83%  '

# Row 44
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').Value = '''63.97'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.88%  '

# Row 45
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').Value = '''5.03'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.85%  '

# Row 46
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '''0.202'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.07%  '

# Row 47
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '''8.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.01%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.101'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.31%  '

# Row 49
$ws.Range('B49').Value = 'BinanceUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D49').Value = '''1.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.10%  '

# Row 50
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').Value = '''1.15'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.74%  '

# Row 51
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = '''4.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.18%  '
